# "hello may not be included"
# Add a new column C to Sheet1 containing sequential row numbers (1-14),
# next to the existing day-of-week / month columns (A/B).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C, rows 1-14: numbers 1 through 14 (one per row).
for ($r = 1; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = $r
}

# Narrow column C to fit the small numbers.
$ws.Columns.Item(3).ColumnWidth = 4.6640625

# Leave the selection where it ends up after filling the column (just
# past the new data, matching the editor's cursor position).
$ws.Range("D15").Select()
